# Applies the table-content edits described by the commit diff:
#  - Table 2 / Table 3 (the duplicated "iris descr" cross-tab): the
#    second header row's group-count cell ("My custom group label")
#    gets the count annotation "(N=50)" filled in.
#  - Table 4 (the categorical "a"/"b" comparison table): the 5th grid
#    column (the "p" column) narrows from 1291 -> 1212 twips, and the
#    "a"/"b" rows' counts / percentages / p-value / CI get updated.

$d = $word.ActiveDocument

# --- Table 2, row 2, cell 2: empty -> "(N=50)" -------------------------
$t2 = $d.Tables.Item(2)
$t2.Rows.Item(2).Cells.Item(2).Range.Text = "(N=50)"

# --- Table 3, row 2, cell 2: empty -> "(N=50)" -------------------------
$t3 = $d.Tables.Item(3)
$t3.Rows.Item(2).Cells.Item(2).Range.Text = "(N=50)"

# --- Table 4: narrow the 5th grid column (1291 -> 1212 twips) ----------
# Column.Width is expressed in points (1 pt = 20 twips), so 1212 twips
# == 60.6 pt.
$t4 = $d.Tables.Item(4)
$t4.Columns.Item(5).Width = 60.6

# --- Table 4, row 18 ("a" row) ------------------------------------------
$rowA = $t4.Rows.Item(18)
$rowA.Cells.Item(2).Range.Text = "17 (57%)"
$rowA.Cells.Item(3).Range.Text = "14 (47%)"
$rowA.Cells.Item(4).Range.Text = "31 (52%)"
$rowA.Cells.Item(5).Range.Text = "0.438"
$rowA.Cells.Item(6).Range.Text = "[-0.15, 0.35]"

# --- Table 4, row 19 ("b" row) ------------------------------------------
$rowB = $t4.Rows.Item(19)
$rowB.Cells.Item(2).Range.Text = "13 (43%)"
$rowB.Cells.Item(3).Range.Text = "16 (53%)"
$rowB.Cells.Item(4).Range.Text = "29 (48%)"
